$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Cell, [string]$Text, [bool]$ForceText)
    if ($ForceText) {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Text
        $Cell.ClearFormats()
    } else {
        $Cell.Value = $Text
    }
}

Set-CellText $ws.Range("D2") "25.195.45" $false
Set-CellText $ws.Range("E2") "  -2.79%  " $false
Set-CellText $ws.Range("D3") "1.656.11" $false
Set-CellText $ws.Range("E3") "  -4.57%  " $false
Set-CellText $ws.Range("D4") "0.9975" $true
Set-CellText $ws.Range("E4") "  -0.16%  " $false
Set-CellText $ws.Range("D5") "234.59" $true
Set-CellText $ws.Range("E5") "  -4.73%  " $false
Set-CellText $ws.Range("D6") "0.9983" $true
Set-CellText $ws.Range("E6") "  -0.13%  " $false
Set-CellText $ws.Range("D7") "0.4783" $true
Set-CellText $ws.Range("E7") "  -4.90%  " $false
Set-CellText $ws.Range("D8") "0.2575" $true
Set-CellText $ws.Range("E8") "  -5.72%  " $false
Set-CellText $ws.Range("D9") "0.06125" $true
Set-CellText $ws.Range("E9") "  -0.71%  " $false
Set-CellText $ws.Range("D10") "0.07058" $true
Set-CellText $ws.Range("E10") "  -2.55%  " $false
Set-CellText $ws.Range("D11") "1.644.29" $false
Set-CellText $ws.Range("E11") "  -5.32%  " $false
Set-CellText $ws.Range("D12") "14.51" $true
Set-CellText $ws.Range("E12") "  -4.46%  " $false
Set-CellText $ws.Range("D13") "0.5793" $true
Set-CellText $ws.Range("E13") "  -11.55%  " $false
Set-CellText $ws.Range("D14") "4.331" $true
Set-CellText $ws.Range("E14") "  -9.21%  " $false
Set-CellText $ws.Range("D15") "73.89" $true
Set-CellText $ws.Range("E15") "  -4.19%  " $false
Set-CellText $ws.Range("D16") "0.9994" $true
Set-CellText $ws.Range("E16") "  +0.09%  " $false
Set-CellText $ws.Range("D17") "0.9988" $true
Set-CellText $ws.Range("E17") "  +0.00%  " $false
Set-CellText $ws.Range("D18") "25.155.32" $false
Set-CellText $ws.Range("E18") "  -3.01%  " $false
Set-CellText $ws.Range("D19") "0.000006657" $true
Set-CellText $ws.Range("E19") "  -2.37%  " $false
Set-CellText $ws.Range("D20") "11.33" $true
Set-CellText $ws.Range("E20") "  -4.53%  " $false
Set-CellText $ws.Range("D21") "1.853.95" $false
Set-CellText $ws.Range("E21") "  -5.50%  " $false
Set-CellText $ws.Range("D22") "4.345" $true
Set-CellText $ws.Range("E22") "  -5.30%  " $false
Set-CellText $ws.Range("D23") "8.521" $true
Set-CellText $ws.Range("E23") "  -3.08%  " $false
Set-CellText $ws.Range("D24") "5.266" $true
Set-CellText $ws.Range("E24") "  -3.75%  " $false
Set-CellText $ws.Range("D25") "134.62" $true
Set-CellText $ws.Range("E25") "  +0.59%  " $false
Set-CellText $ws.Range("D26") "15.03" $true
Set-CellText $ws.Range("E26") "  -1.26%  " $false
Set-CellText $ws.Range("D27") "1.376" $true
Set-CellText $ws.Range("E27") "  -4.58%  " $false
Set-CellText $ws.Range("D28") "104.35" $true
Set-CellText $ws.Range("E28") "  -0.88%  " $false
Set-CellText $ws.Range("D29") "1.656" $true
Set-CellText $ws.Range("E29") "  -7.44%  " $false
Set-CellText $ws.Range("D30") "3.929" $true
Set-CellText $ws.Range("E30") "  -1.56%  " $false
Set-CellText $ws.Range("D31") "0.07621" $true
Set-CellText $ws.Range("E31") "  -6.09%  " $false
Set-CellText $ws.Range("D32") "3.569" $true
Set-CellText $ws.Range("E32") "  -3.49%  " $false
Set-CellText $ws.Range("D33") "0.9985" $true
Set-CellText $ws.Range("E33") "  +0.04%  " $false
Set-CellText $ws.Range("D34") "0.04317" $true
Set-CellText $ws.Range("E34") "  -8.75%  " $false
Set-CellText $ws.Range("D35") "2.597" $true
Set-CellText $ws.Range("E35") "  -2.15%  " $false
Set-CellText $ws.Range("D36") "0.9393" $true
Set-CellText $ws.Range("E36") "  -5.91%  " $false
Set-CellText $ws.Range("D37") "0.5983" $true
Set-CellText $ws.Range("E37") "  -2.31%  " $false
Set-CellText $ws.Range("D38") "2.587" $true
Set-CellText $ws.Range("E38") "  -5.98%  " $false
Set-CellText $ws.Range("D39") "0.8554" $true
Set-CellText $ws.Range("E39") "  -3.44%  " $false
Set-CellText $ws.Range("D40") "0.9993" $true
Set-CellText $ws.Range("E40") "  +0.09%  " $false
Set-CellText $ws.Range("D41") "0.01485" $true
Set-CellText $ws.Range("E41") "  -7.39%  " $false
Set-CellText $ws.Range("D42") "98.99" $true
Set-CellText $ws.Range("E42") "  -1.66%  " $false
Set-CellText $ws.Range("D43") "1.799" $true
Set-CellText $ws.Range("E43") "  -8.27%  " $false
Set-CellText $ws.Range("D44") "0.3704" $true
Set-CellText $ws.Range("E44") "  -5.34%  " $false
Set-CellText $ws.Range("D45") "4.644" $true
Set-CellText $ws.Range("E45") "  -7.34%  " $false
Set-CellText $ws.Range("B46") "Aptos" $false
Set-CellText $ws.Range("C46") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" $false
Set-CellText $ws.Range("D46") "6.142" $true
Set-CellText $ws.Range("E46") "  -2.93%  " $false
Set-CellText $ws.Range("B47") "Algorand" $false
Set-CellText $ws.Range("C47") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" $false
Set-CellText $ws.Range("D47") "0.1100" $true
Set-CellText $ws.Range("E47") "  -6.69%  " $false
Set-CellText $ws.Range("D48") "0.05232" $true
Set-CellText $ws.Range("E48") "  -0.92%  " $false
Set-CellText $ws.Range("D49") "29.22" $true
Set-CellText $ws.Range("E49") "  -4.96%  " $false
Set-CellText $ws.Range("B50") "NEARProtocol" $false
Set-CellText $ws.Range("C50") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" $false
Set-CellText $ws.Range("D50") "1.210" $true
Set-CellText $ws.Range("E50") "  -1.92%  " $false
Set-CellText $ws.Range("B51") "TrueUSD" $false
Set-CellText $ws.Range("C51") "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd" $false
Set-CellText $ws.Range("E51") "  -0.01%  " $false
